$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.687.55'
$ws.Range("E2").Value = '  -3.55%  '

$ws.Range("D3").Value = '1.743.17'
$ws.Range("E3").Value = '  -5.50%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '237.50'
$ws.Range("E5").Value = '  -8.61%  '

$ws.Range("D7").Value = '0.4929'
$ws.Range("E7").Value = '  -6.65%  '

$ws.Range("D8").Value = '41.59'
$ws.Range("E8").Value = '  -7.58%  '

$ws.Range("D9").Value = '0.2376'
$ws.Range("E9").Value = '  -24.75%  '

$ws.Range("D10").Value = '0.05957'
$ws.Range("E10").Value = '  -12.39%  '

$ws.Range("D11").Value = '1.740.43'
$ws.Range("E11").Value = '  -5.82%  '

$ws.Range("E12").Value = '  -12.22%  '

$ws.Range("D13").Value = '14.59'
$ws.Range("E13").Value = '  -23.02%  '

$ws.Range("D14").Value = '4.454'
$ws.Range("E14").Value = '  -11.19%  '

$ws.Range("D15").Value = '77.08'
$ws.Range("E15").Value = '  -12.71%  '

$ws.Range("D16").Value = '0.5730'
$ws.Range("E16").Value = '  -27.05%  '

$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  +0.07%  '

$ws.Range("E18").Value = '  +0.12%  '

$ws.Range("D19").Value = '25.733.76'
$ws.Range("E19").Value = '  -3.43%  '

$ws.Range("D20").Value = '11.44'
$ws.Range("E20").Value = '  -17.62%  '

$ws.Range("D21").Value = '0.000006437'
$ws.Range("E21").Value = '  -18.75%  '

$ws.Range("D22").Value = '1.959.90'
$ws.Range("E22").Value = '  -6.28%  '

$ws.Range("D23").Value = '3.952'
$ws.Range("E23").Value = '  -14.20%  '

$ws.Range("D24").Value = '5.041'
$ws.Range("E24").Value = '  -15.77%  '

$ws.Range("D25").Value = '7.760'
$ws.Range("E25").Value = '  -16.99%  '

$ws.Range("D26").Value = '136.58'
$ws.Range("E26").Value = '  -4.50%  '

$ws.Range("D27").Value = '1.475'
$ws.Range("E27").Value = '  -12.06%  '

$ws.Range("D28").Value = '1.827'
$ws.Range("E28").Value = '  -18.00%  '

$ws.Range("D29").Value = '14.51'
$ws.Range("E29").Value = '  -14.84%  '

$ws.Range("D30").Value = '100.63'
$ws.Range("E30").Value = '  -9.37%  '

$ws.Range("D31").Value = '3.781'
$ws.Range("E31").Value = '  -10.06%  '

$ws.Range("D32").Value = '0.08137'
$ws.Range("E32").Value = '  -6.63%  '

$ws.Range("D33").Value = '3.340'

$ws.Range("D34").Value = '0.04354'
$ws.Range("E34").Value = '  -10.86%  '

$ws.Range("D35").Value = '1.000'
$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").Value = '2.679'
$ws.Range("E36").Value = '  -6.38%  '

$ws.Range("D37").Value = '1.016'
$ws.Range("E37").Value = '  -11.00%  '

$ws.Range("D38").Value = '0.6062'
$ws.Range("E38").Value = '  -17.16%  '

$ws.Range("D39").Value = '2.714'
$ws.Range("E39").Value = '  -12.68%  '

$ws.Range("D40").Value = '2.078'
$ws.Range("E40").Value = '  -9.48%  '

$ws.Range("D41").Value = '1.002'
$ws.Range("E41").Value = '  +0.11%  '

$ws.Range("D42").Value = '103.20'
$ws.Range("E42").Value = '  -6.12%  '

$ws.Range("D43").Value = '0.01480'
$ws.Range("E43").Value = '  -14.54%  '

$ws.Range("D44").Value = '0.7782'
$ws.Range("E44").Value = '  -13.68%  '

$ws.Range("D45").Value = '5.137'
$ws.Range("E45").Value = '  -13.50%  '

$ws.Range("D46").Value = '0.3766'
$ws.Range("E46").Value = '  -21.81%  '

$ws.Range("D47").Value = '0.05108'
$ws.Range("E47").Value = '  -12.31%  '

$ws.Range("D48").Value = '5.968'
$ws.Range("E48").Value = '  -22.60%  '

$ws.Range("D49").Value = '0.1068'
$ws.Range("E49").Value = '  -14.06%  '

$ws.Range("D50").Value = '30.28'
$ws.Range("E50").Value = '  -13.04%  '

$ws.Range("D51").Value = '52.57'
$ws.Range("E51").Value = '  -12.35%  '
